$d = $word.ActiveDocument

# Update the zip changeset number: C37714 -> C37722
# The text is split across 4 bold runs: "C" "3" "771" "4".
# We need "771" -> "77" and "4" -> "22", each staying in its own run
# (not merged with neighboring runs that share identical formatting).
#
# Directly assigning Range.Text causes the engine to coalesce the edited
# run with adjacent runs that have identical formatting. To avoid that,
# we briefly toggle Bold off on the target run before editing its text,
# then turn Bold back on for the new text -- this keeps the run boundary
# intact because at edit-time the run's formatting differs from its
# neighbors.

$r = $d.Content
$found = $r.Find.Execute("771")
if ($found) {
    $start1 = $r.Start
    $end1 = $r.End

    # "771" -> "77"
    $r1 = $d.Range($start1, $end1)
    $r1.Font.Bold = $false
    $r1.Text = "77"
    $newEnd1 = $start1 + 2
    $r1b = $d.Range($start1, $newEnd1)
    $r1b.Font.Bold = $true

    # the following run contains "4" (part of the same "C3771" + "4" = "C37714"
    # changeset number); change it to "22" -> "C37722"
    $r4 = $d.Range($newEnd1, $newEnd1 + 1)
    $r4.Font.Bold = $false
    $r4.Text = "22"
    $r4b = $d.Range($newEnd1, $newEnd1 + 2)
    $r4b.Font.Bold = $true
}
